$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Bata no liquidificador as cenouras, os ovos e o óleo." -> replace the
#    trailing "." with " até obter uma mistura homogênea. "
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(11)
$p1End = $p1.Range.End
$dotRange = $d.Range($p1End - 2, $p1End - 1)
$dotRange.Text = " até obter uma mistura homogênea. "

# ---------------------------------------------------------------------------
# 2) "Em uma tigela, misture o açúcar e a farinha." -> replace the trailing
#    "." with " de trigo."
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(12)
$p2End = $p2.Range.End
$dotRange2 = $d.Range($p2End - 2, $p2End - 1)
$dotRange2.Text = " de trigo."

# ---------------------------------------------------------------------------
# 3) Add four new numbered steps after "Em uma tigela..." paragraph
# ---------------------------------------------------------------------------
$lastList = $d.Paragraphs.Item(12)
$r = $lastList.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$d.Paragraphs.Item(13).Range.Text = "Adicione a mistura liquida á tigela e mexa bem."

$lastList = $d.Paragraphs.Item(13)
$r = $lastList.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$d.Paragraphs.Item(14).Range.Text = "Acrescente fermento e misture delicadamente."

$lastList = $d.Paragraphs.Item(14)
$r = $lastList.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$d.Paragraphs.Item(15).Range.Text = "Despeje a massa em uma forma untada e enfarinhada."

$lastList = $d.Paragraphs.Item(15)
$r = $lastList.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$d.Paragraphs.Item(16).Range.Text = "Asse em forno preaquecido a 180°C por cerca de 40 minutos."

# ---------------------------------------------------------------------------
# 4) Add a blank (non-numbered) paragraph after the last step
# ---------------------------------------------------------------------------
$lastList = $d.Paragraphs.Item(16)
$r = $lastList.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$d.Paragraphs.Item(17).Range.ListFormat.RemoveNumbers()

# ---------------------------------------------------------------------------
# 5) Add the "OBSERVAÇÃO:" paragraph (bold label + normal remark text)
# ---------------------------------------------------------------------------
$blank = $d.Paragraphs.Item(17)
$r = $blank.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$obsPara = $d.Paragraphs.Item(18)
$obsPara.Range.ListFormat.RemoveNumbers()
$obsPara.Range.Text = "OBSERVAÇÃO:"
$obsPara.Range.Font.Bold = 1
$tailRange = $d.Range($obsPara.Range.End - 1, $obsPara.Range.End - 1)
$tailRange.InsertAfter(" Adicione uma cobertura simples de chocolate após o bolo esfriar. ")
$tailRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# 6) Add the final trailing blank paragraph
# ---------------------------------------------------------------------------
$obsPara2 = $d.Paragraphs.Item(18)
$r = $obsPara2.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$d.Paragraphs.Item(19).Range.ListFormat.RemoveNumbers()

$d.Paragraphs.Item(19).Range.Text
